# Apply cryptocurrency price/volume updates as described in the commit
# "Updated cryptos list on Mon Jul 17 09:55:31 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.323.14"
$ws.Range("E2").Value = "'  +0.06%  "

$ws.Range("D3").Value = "'1.928.54"
$ws.Range("E3").Value = "'  -0.07%  "

$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "'  +0.54%  "

$ws.Range("D5").Value = "'0.7348"
$ws.Range("E5").Value = "'  +2.13%  "

$ws.Range("D6").Value = "'243.57"
$ws.Range("E6").Value = "'  -2.45%  "

$ws.Range("D7").Value = "'1.004"
$ws.Range("E7").Value = "'  +0.41%  "

$ws.Range("D8").Value = "'27.52"
$ws.Range("E8").Value = "'  -1.47%  "

$ws.Range("D9").Value = "'0.3144"
$ws.Range("E9").Value = "'  -1.82%  "

$ws.Range("D10").Value = "'0.07003"
$ws.Range("E10").Value = "'  -1.42%  "

$ws.Range("D11").Value = "'0.08034"
$ws.Range("E11").Value = "'  +0.16%  "

$ws.Range("D12").Value = "'0.7731"
$ws.Range("E12").Value = "'  -1.92%  "

$ws.Range("D13").Value = "'1.977.18"
$ws.Range("E13").Value = "'  +2.43%  "

$ws.Range("D14").Value = "'5.352"
$ws.Range("E14").Value = "'  -0.49%  "

$ws.Range("D15").Value = "'92.88"
$ws.Range("E15").Value = "'  -1.92%  "

$ws.Range("D16").Value = "'14.40"
$ws.Range("E16").Value = "'  -1.69%  "

$ws.Range("D17").Value = "'30.312.75"
$ws.Range("E17").Value = "'  +0.04%  "

$ws.Range("D18").Value = "'249.69"
$ws.Range("E18").Value = "'  -2.69%  "

$ws.Range("D19").Value = "'5.930"
$ws.Range("E19").Value = "'  +3.33%  "

$ws.Range("D20").Value = "'0.000007910"
$ws.Range("E20").Value = "'  -1.98%  "

$ws.Range("B21").Value = "'Dai"
$ws.Range("C21").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'1.004"
$ws.Range("E21").Value = "'  +0.35%  "

$ws.Range("B22").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "'2.139.24"
$ws.Range("E22").Value = "'  -1.90%  "

$ws.Range("D23").Value = "'1.006"
$ws.Range("E23").Value = "'  +0.50%  "

$ws.Range("D24").Value = "'6.625"
$ws.Range("E24").Value = "'  -2.81%  "

$ws.Range("D25").Value = "'9.530"
$ws.Range("E25").Value = "'  -0.26%  "

$ws.Range("D26").Value = "'165.65"
$ws.Range("E26").Value = "'  +0.76%  "

$ws.Range("D27").Value = "'18.98"
$ws.Range("E27").Value = "'  -0.66%  "

$ws.Range("D28").Value = "'0.1275"
$ws.Range("E28").Value = "'  -0.70%  "

$ws.Range("D29").Value = "'2.160"
$ws.Range("E29").Value = "'  -6.26%  "

$ws.Range("D30").Value = "'1.359"
$ws.Range("E30").Value = "'  +0.41%  "

$ws.Range("D31").Value = "'1.557"
$ws.Range("E31").Value = "'  +1.56%  "

$ws.Range("E32").Value = "'  -0.38%  "

$ws.Range("D33").Value = "'4.095"
$ws.Range("E33").Value = "'  -1.43%  "

$ws.Range("D34").Value = "'0.05200"
$ws.Range("E34").Value = "'  +1.82%  "

$ws.Range("D35").Value = "'1.303"
$ws.Range("E35").Value = "'  +0.86%  "

$ws.Range("D36").Value = "'0.7515"
$ws.Range("E36").Value = "'  +0.22%  "

$ws.Range("D37").Value = "'2.771"
$ws.Range("E37").Value = "'  +0.08%  "

$ws.Range("D38").Value = "'0.01943"
$ws.Range("E38").Value = "'  -2.19%  "

$ws.Range("D39").Value = "'2.789"

$ws.Range("D40").Value = "'6.507"
$ws.Range("E40").Value = "'  +1.70%  "

$ws.Range("D41").Value = "'76.36"
$ws.Range("E41").Value = "'  -2.61%  "

$ws.Range("D42").Value = "'0.4483"
$ws.Range("E42").Value = "'  -0.93%  "

$ws.Range("D43").Value = "'1.945"
$ws.Range("E43").Value = "'  -2.54%  "

$ws.Range("D44").Value = "'1.003"
$ws.Range("E44").Value = "'  +0.34%  "

$ws.Range("D45").Value = "'0.8386"
$ws.Range("E45").Value = "'  -0.84%  "

$ws.Range("D46").Value = "'7.654"
$ws.Range("E46").Value = "'  +2.23%  "

$ws.Range("D47").Value = "'9.946"
$ws.Range("E47").Value = "'  +1.20%  "

$ws.Range("E48").Value = "'  -0.23%  "

$ws.Range("D49").Value = "'37.41"
$ws.Range("E49").Value = "'  +1.53%  "

$ws.Range("D50").Value = "'2.074.71"
$ws.Range("E50").Value = "'  -1.07%  "

$ws.Range("D51").Value = "'0.1219"
$ws.Range("E51").Value = "'  +6.47%  "

Write-Output "Updated cryptos list"